# Change font filename extensions from .otf to .ttf in column B (Font column)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*.otf") {
        $newVal = $val -replace "\.otf$", ".ttf"
        $cell.Value2 = $newVal
    }
}

# Scroll the view so that row 74 is near the top, matching the author's saved
# view state (best effort - some COM hosts do not persist pure scroll
# position without an accompanying freeze/split pane).
try {
    $excel.ActiveWindow.ScrollRow = 74
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
